$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.624.35'
$ws.Range("E2").Value = '  -2.21%  '
$ws.Range("D3").Value = '1.793.84'
$ws.Range("E3").Value = '  -2.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.76%  '
$ws.Range("E6").Value = '  -2.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2768'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06746'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07525'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.69%  '
$ws.Range("D12").Value = '1.794.51'
$ws.Range("E12").Value = '  -2.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.799'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6140'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.41%  '
$ws.Range("D15").Value = '2.036.35'
$ws.Range("E15").Value = '  -2.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '75.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009032'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.59%  '
$ws.Range("D18").Value = '28.598.57'
$ws.Range("E18").Value = '  -2.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.430'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '210.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.807'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.004'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.077'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1259'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.416'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06115'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.78%  '
$ws.Range("E31").Value = '  -1.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.811'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.786'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.735'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.048'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6414'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.81%  '
$ws.Range("E37").Value = '  -1.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.711'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.420'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01694'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.09%  '
$ws.Range("D41").Value = '1.141.18'
$ws.Range("E41").Value = '  -6.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8809'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.006'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = '1.944.14'
$ws.Range("E45").Value = '  -2.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000111'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.584'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05487'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.342'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4477'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.73%  '
